# Weekly update of Fruta/Hortaliza (Guayaba) data: dates and prices for several rows
# have been shuffled to reflect a new weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44405

# Row 4
$ws.Range("D4").Value = 44431
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 1300
$ws.Range("O4").Value = 1300
$ws.Range("P4").Value = 1300
$ws.Range("S4").Value = 1300

# Row 5
$ws.Range("D5").Value = 44476
$ws.Range("M5").Value = 80

# Row 6
$ws.Range("D6").Value = 44438
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44417
$ws.Range("M7").Value = 80

# Row 8
$ws.Range("D8").Value = 44432
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 1300
$ws.Range("O8").Value = 1300
$ws.Range("P8").Value = 1300
$ws.Range("S8").Value = 1300

# Row 11
$ws.Range("D11").Value = 44473
$ws.Range("M11").Value = 120

# Row 12
$ws.Range("D12").Value = 44418
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 1200
$ws.Range("O12").Value = 1200
$ws.Range("P12").Value = 1200
$ws.Range("S12").Value = 1200

# Row 13
$ws.Range("D13").Value = 44424
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 1200
$ws.Range("O13").Value = 1200
$ws.Range("P13").Value = 1200
$ws.Range("S13").Value = 1200
